$d = $word.ActiveDocument

# Fill in the Name / ID / Class table cells (second column) which were empty.
$table = $d.Tables.Item(1)
$table.Cell(1, 2).Range.Text = "Chua Zhi Yang"
$table.Cell(2, 2).Range.Text = "2004726B"
$table.Cell(3, 2).Range.Text = "P01"

# Remove the leftover "_GoBack" bookmark (Word drops this automatically on save).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
